$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.552136229771481
$ws.Range("C2").Value = 1.493186813606976
$ws.Range("D2").Value = 4.559132871428442
$ws.Range("E2").Value = 1.456920224378268
$ws.Range("F2").Value = 1.550484622924876
$ws.Range("G2").Value = 1.570172298926858
$ws.Range("H2").Value = 1.457347192194637
$ws.Range("B3").Value = 1.563193023149681
$ws.Range("C3").Value = 1.500165029430642
$ws.Range("D3").Value = 3.651525841603914
$ws.Range("E3").Value = 1.458610037681003
$ws.Range("F3").Value = 1.55689358682896
$ws.Range("G3").Value = 1.558931202156377
$ws.Range("H3").Value = 1.458595287136383
$ws.Range("B4").Value = 1.526684029639687
$ws.Range("C4").Value = 1.549274810175597
$ws.Range("D4").Value = 4.19838724095648
$ws.Range("E4").Value = 1.451047547992394
$ws.Range("F4").Value = 1.527468092681718
$ws.Range("G4").Value = 1.57236504845805
$ws.Range("H4").Value = 1.451431135937862
$ws.Range("B5").Value = 1.545217619614819
$ws.Range("C5").Value = 1.601553576152952
$ws.Range("D5").Value = 2.42626073629948
$ws.Range("E5").Value = 1.446216490933229
$ws.Range("F5").Value = 1.537805340834631
$ws.Range("G5").Value = 1.57880164550358
$ws.Range("H5").Value = 1.446166112413758
$ws.Range("B6").Value = 1.551005964670673
$ws.Range("C6").Value = 1.602459148828121
$ws.Range("D6").Value = 0.6342104514815061
$ws.Range("E6").Value = 1.448875947771873
$ws.Range("F6").Value = 1.543419494542331
$ws.Range("G6").Value = 1.548287313870727
$ws.Range("H6").Value = 1.448580771627436
$ws.Range("B7").Value = 1.549486314736102
$ws.Range("C7").Value = 1.641837099424327
$ws.Range("D7").Value = 1.000518488509308
$ws.Range("E7").Value = 1.461064592524589
$ws.Range("F7").Value = 1.54268538216686
$ws.Range("G7").Value = 1.587202116579776
$ws.Range("H7").Value = 1.460838994595111
$ws.Range("B8").Value = 1.573590969358645
$ws.Range("C8").Value = 1.25963316010294
$ws.Range("D8").Value = 1.398601356383271
$ws.Range("E8").Value = 1.463279494679763
$ws.Range("F8").Value = 1.565222802774
$ws.Range("G8").Value = 1.347264377506427
$ws.Range("H8").Value = 1.463289021859106
$ws.Range("B9").Value = 1.563634562990191
$ws.Range("C9").Value = 1.549120835106842
$ws.Range("D9").Value = 1.032885358940224
$ws.Range("E9").Value = 1.455019439359377
$ws.Range("F9").Value = 1.557877905496985
$ws.Range("G9").Value = 1.524252566768621
$ws.Range("H9").Value = 1.455487230471159
$ws.Range("B10").Value = 1.344117272470109
$ws.Range("C10").Value = 1.639123914208648
$ws.Range("D10").Value = 2.919303411800664
$ws.Range("E10").Value = 1.44357601379494
$ws.Range("F10").Value = 1.375272804597484
$ws.Range("G10").Value = 1.602662455506086
$ws.Range("H10").Value = 1.444906915254569
$ws.Range("B11").Value = 1.298091163975458
$ws.Range("C11").Value = 1.646314611292481
$ws.Range("D11").Value = 2.465559063346524
$ws.Range("E11").Value = 1.446995113239449
$ws.Range("F11").Value = 1.329380741006684
$ws.Range("G11").Value = 1.599602480699534
$ws.Range("H11").Value = 1.447518338027816
$ws.Range("B12").Value = 1.073177682413108
$ws.Range("C12").Value = 1.6308334743146
$ws.Range("D12").Value = 0.9315824148048609
$ws.Range("E12").Value = 1.42465114317971
$ws.Range("F12").Value = 1.082554188446728
$ws.Range("G12").Value = 1.56334331424269
$ws.Range("H12").Value = 1.424266206036189
$ws.Range("B13").Value = 1.329072515563819
$ws.Range("C13").Value = 1.636853583344017
$ws.Range("D13").Value = 2.301587582174857
$ws.Range("E13").Value = 1.440235713151947
$ws.Range("F13").Value = 1.350229857885922
$ws.Range("G13").Value = 1.587080947628144
$ws.Range("H13").Value = 1.440926698548765
